$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2 value to the new part reference (adds a new shared string entry)
$ws.Range("D2").Value = "C1509219"

# Move the active selection to D3 (as reflected in the diff)
$ws.Range("D3").Select()
